{"js": "// The document opens with a masthead made of two paragraphs:\n//   1) \"On Pilgrimage,  \\n October ==================\"  (italic \"On Pilgrimage\")\n//   2) \"By Dorothy Day\"  (bold)\n//\n// The edit collapses that into a pandoc-style title block:\n//   1) \"October\"        -- paragraph styled as \"Title\"\n//   2) \"% Dorothy Day\"  -- plain paragraph (pandoc author line)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst authorPara = paragraphs.items[1];\n\n// Replace the whole first paragraph's content with \"October\" and apply\n// the built-in \"Title\" style, wiping out the previous run-level formatting\n// (italic \"On Pilgrimage\", the manual line break, etc.).\nconst titleRange = titlePara.getRange();\ntitleRange.clear();\ntitleRange.insertText(\"October\", Word.InsertLocation.replace);\ntitlePara.style = \"Title\";\n\n// Replace the second paragraph's content with the pandoc author line,\n// dropping the previous bold formatting.\nconst authorRange = authorPara.getRange();\nauthorRange.clear();\nauthorRange.insertText(\"% Dorothy Day\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document opens with a masthead made of two paragraphs:\n#   1) \"On Pilgrimage,  \" + manual line break + \" October ==================\"\n#      (with \"On Pilgrimage\" in italics)\n#   2) \"By Dorothy Day\"  (bold)\n#\n# Turn that into a pandoc-style title block:\n#   1) \"October\"        -- paragraph styled as \"Title\"\n#   2) \"% Dorothy Day\"  -- plain paragraph (pandoc author line)\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: masthead -> Title-styled \"October\" ---\n$titlePara = $d.Paragraphs(1)\n$titleRange = $titlePara.Range\n# Exclude the trailing paragraph mark so we only touch the paragraph's text.\n$titleText = $d.Range($titleRange.Start, $titleRange.End - 1)\n$titleText.Delete()\n$titleText.InsertBefore(\"October\")\n$titlePara.Style = \"Title\"\n\n# --- Paragraph 2: \"By Dorothy Day\" -> \"% Dorothy Day\" ---\n$authorPara = $d.Paragraphs(2)\n$authorRange = $authorPara.Range\n$authorText = $d.Range($authorRange.Start, $authorRange.End - 1)\n$authorText.Delete()\n$authorText.InsertBefore(\"% Dorothy Day\")\n"}
